$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" "59.814.96"
Set-TextValue $ws "E2" "  +2.59%  "

Set-TextValue $ws "D3" "2.419.23"
Set-TextValue $ws "E3" "  +2.69%  "

Set-TextValue $ws "E4" "  +0.06%  "

Set-TextValue $ws "D5" "553.67"
Set-TextValue $ws "E5" "  +2.31%  "

Set-TextValue $ws "D6" "138.03"

Set-TextValue $ws "E7" "  -0.09%  "

Set-TextValue $ws "D8" "0.568"
Set-TextValue $ws "E8" "  +1.02%  "

Set-TextValue $ws "D9" "0.106"
Set-TextValue $ws "E9" "  +3.91%  "

Set-TextValue $ws "D10" "5.82"
Set-TextValue $ws "E10" "  +3.79%  "

Set-TextValue $ws "D11" "0.359"
Set-TextValue $ws "E11" "  +1.44%  "

Set-TextValue $ws "E12" "  -2.18%  "

Set-TextValue $ws "D13" "24.67"

Set-TextValue $ws "D14" "2.852.33"
Set-TextValue $ws "E14" "  +2.77%  "

Set-TextValue $ws "D15" "59.743.28"
Set-TextValue $ws "E15" "  +2.52%  "

Set-TextValue $ws "E16" "  +4.11%  "

Set-TextValue $ws "D17" "2.409.55"
Set-TextValue $ws "E17" "  +1.57%  "

Set-TextValue $ws "D18" "11.36"
Set-TextValue $ws "E18" "  +5.82%  "

Set-TextValue $ws "D19" "4.43"
Set-TextValue $ws "E19" "  +3.35%  "

Set-TextValue $ws "D20" "333.56"
Set-TextValue $ws "E20" "  +0.16%  "

Set-TextValue $ws "D21" "6.93"
Set-TextValue $ws "E21" "  +1.87%  "

Set-TextValue $ws "E22" "  -0.03%  "

Set-TextValue $ws "D23" "64.63"
Set-TextValue $ws "E23" "  +2.40%  "

Set-TextValue $ws "E24" "  +1.40%  "

Set-TextValue $ws "D25" "8.60"
Set-TextValue $ws "E25" "  +1.41%  "

Set-TextValue $ws "E26" "  +0.04%  "

Set-TextValue $ws "D27" "1.37"
Set-TextValue $ws "E27" "  -1.74%  "

$subscriptThree = [string]([char]0x2083)
Set-TextValue $ws "D28" ("0.0" + $subscriptThree + "0785")
Set-TextValue $ws "E28" "  +6.54%  "

Set-TextValue $ws "E29" "  +3.24%  "

Set-TextValue $ws "D30" "170.77"
Set-TextValue $ws "E30" "  -0.25%  "

Set-TextValue $ws "D31" "6.27"
Set-TextValue $ws "E31" "  +2.00%  "

Set-TextValue $ws "D32" "18.65"
Set-TextValue $ws "E32" "  +1.05%  "

Set-TextValue $ws "E33" "  +0.02%  "

Set-TextValue $ws "E35" "  +5.50%  "

Set-TextValue $ws "E36" "  -0.66%  "

Set-TextValue $ws "E37" "  +0.20%  "

Set-TextValue $ws "E38" "  -1.18%  "

Set-TextValue $ws "D39" "40.09"
Set-TextValue $ws "E39" "  +2.33%  "

Set-TextValue $ws "D40" "0.423"
Set-TextValue $ws "E40" "  +12.06%  "

Set-TextValue $ws "D41" "313.34"
Set-TextValue $ws "E41" "  +6.56%  "

Set-TextValue $ws "D42" "3.73"
Set-TextValue $ws "E42" "  +2.26%  "

Set-TextValue $ws "D43" "142.62"
Set-TextValue $ws "E43" "  -1.64%  "

Set-TextValue $ws "D44" "0.0963"
Set-TextValue $ws "E44" "  +1.53%  "

Set-TextValue $ws "E45" "  +3.89%  "

Set-TextValue $ws "D46" "0.415"
Set-TextValue $ws "E46" "  +8.66%  "

Set-TextValue $ws "D47" "19.21"
Set-TextValue $ws "E47" "  +0.04%  "

Set-TextValue $ws "E48" "  +1.30%  "

Set-TextValue $ws "E49" "  +2.32%  "

Set-TextValue $ws "E50" "  -0.32%  "

Set-TextValue $ws "D51" "1.61"
Set-TextValue $ws "E51" "  +3.92%  "
